$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 120.7114695641766
    3  = 111.8753144741339
    4  = 111.8330232313298
    5  = 103.0721123017533
    6  = 73.50003484837279
    7  = 73.14110564301146
    8  = 72.53603635816313
    9  = 72.4053588385395
    10 = 72.02062024955674
    11 = 71.0543292599872
    12 = 69.86281313742046
    13 = 69.63302543963597
    14 = 69.49932754219512
    15 = 69.3879728982762
    16 = 69.27705552028992
    17 = 69.24798997526473
    18 = 69.1952601805894
    19 = 69.03717468513139
    20 = 68.82575611034714
    21 = 67.8633445773653
    22 = 66.95791352222486
    23 = 66.95410490560782
    24 = 66.89328678689945
    25 = 66.81745182809752
    26 = 66.43657163097996
    27 = 66.34741950883222
    28 = 66.2756368002285
    29 = 66.24259064836983
    30 = 66.18571242107234
    31 = 66.17210229661404
    32 = 66.00165190272328
    33 = 65.93990681471828
    34 = 65.92609419424426
    35 = 65.40998614189506
    36 = 65.04610289196461
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}
